$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 321. This shifts the existing
# rows 321-393 down to 322-394 (their data moves down by one row, matching
# the new weekly observation being prepended to this block), and leaves a
# blank row 321 to be filled with the newest observation.
$ws.Rows(321).Insert()

# Populate the new row 321 with the latest weekly observation.
$ws.Cells.Item(321, 1).Value = 3
$ws.Cells.Item(321, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(321, 3).Value = "Coquimbo"
$ws.Cells.Item(321, 4).Value = 44798
$ws.Cells.Item(321, 5).Value = 5
$ws.Cells.Item(321, 6).Value = 100112012
$ws.Cells.Item(321, 7).Value = "Espinaca"
$ws.Cells.Item(321, 8).Value = "Sin especificar"
$ws.Cells.Item(321, 9).Value = "Primera"
$ws.Cells.Item(321, 10).Value = 230
$ws.Cells.Item(321, 11).Value = 4000
$ws.Cells.Item(321, 12).Value = 4200
$ws.Cells.Item(321, 13).Value = 4096
$ws.Cells.Item(321, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(321, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(321, 16).Value = 1365
$ws.Cells.Item(321, 17).Value = 3
$ws.Cells.Item(321, 18).Value = "Hortaliza"
